$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing hash values in D2 and D3
$ws.Range("D2").Value = "a84876e7e0593b0995b09045b34c582c3c6bafc096fb444939e8534553babc53"
$ws.Range("D3").Value = "0e08fd8eee36e999283f9dd25c3209735b91da1642e3b67aadf38a56da7ed5d6"

# Row 4 now becomes a new entry (index 3) instead of the old "accept by lobe owner" JSON row
$ws.Range("B4").Value = "23.03.2023 10:08 (CET)"
$ws.Range("C4").Value = "https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/process_model_version/-/commit/69c451ed3cc43181ccf7143ba7e6ada345d38d6f"
$ws.Range("D4").Value = "7d94edb76e6b01e30e2dea18580509d0f3efcd866bfcad2a4b298e040d1da0f6"

# New row 5 (index 4)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "24.03.2023 16:20 (CET)"
$ws.Range("C5").Value = "https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/process_model_version/-/commit/7d95b41a9b432c782ac3bd45ad803e677399ab28"
$ws.Range("D5").Value = "b4ab941a87685341ec282f6e0ff6634ae8a1afa78e3fffdcbad16e151dde0c33"

# New row 6 (index 5) - the old row 4's timestamp/hash data shifted here, with C6 simplified to a plain URL
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "31.03.2023 12:45 (CET)"
$ws.Range("C6").Value = "https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/process_model_version/-/commit/b7ffbf458b88d4980e1c5ee0f583ccc127d4f792"
$ws.Range("D6").Value = "ee5338c88b6ae24c0e57aa6023372a73e87041c67f6f82013b2f688772015316"
